$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Compliance")
$ws2 = $wb.Worksheets.Item("Worksheet")

# --- Compliance sheet: add a new Table1 row (row 20) with the latest test run ---
$lo = $ws1.ListObjects.Item("Table1")
$lo.ListRows.Add() | Out-Null

# Copy formatting down from the previous row so the new row matches styles
$ws1.Range("A19:E19").Copy()
$ws1.Range("A20:E20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A20").Value = 45118.736805555556
$ws1.Range("B20").Value = 1446235
$ws1.Range("C20").Value = 8326
$ws1.Range("D20").Formula = "=Table1[[#This Row],[Failures     ]]/Table1[[#This Row],[Tests     ]]"
$ws1.Range("E20").Formula = "=(`$C`$7-Table1[[#This Row],[Failures     ]])/`$C`$7"

# --- Worksheet sheet: update the DAA reference test entry (was the SCF/CCF entry) ---
$ws2.Range("C2").Value = "4121FA09601D59A55B8D7990009A9D29"
$ws2.Range("C4").Value = "4121FA09601D59A55B8D799094A09D29"
$ws2.Range("C3").Value = "4121FA09601D59A55B8D799055009D29"
$ws2.Range("D3").Value = "'55"
$ws2.Range("D4").Value = "'94"
$ws2.Range("C6").Value = "A: 9A -> 00"
$ws2.Range("E3").Value = "_ Z _ H _ P _ C"
$ws2.Range("E4").Value = "S _ _ H _ P _ _"

# --- Update selections to match the saved state ---
$ws2.Range("A7").Select()
$ws1.Range("A21").Select()
